$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("G3","H3","D4","E4","D5","E5","D6","E6","D7","E7","H8","D9","E9","H10","D11","E11","H12","D13","E13","H14","H15","H16","D17","E17","H18")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = 1
}
